# Insert a new data row at row 581 (a new weekly "Betarraga" price entry),
# pushing the existing rows 581-675 down to 582-676 (dimension grows to R676).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(581).Insert()

$ws.Cells.Item(581, 1).Value = 5
$ws.Cells.Item(581, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(581, 3).Value = "Maule"
$ws.Cells.Item(581, 4).Value = 45218
$ws.Cells.Item(581, 5).Value = 7
$ws.Cells.Item(581, 6).Value = 100114014
$ws.Cells.Item(581, 7).Value = "Betarraga"
$ws.Cells.Item(581, 8).Value = "Sin especificar"
$ws.Cells.Item(581, 9).Value = "Primera"
$ws.Cells.Item(581, 10).Value = 4000
$ws.Cells.Item(581, 11).Value = 500
$ws.Cells.Item(581, 12).Value = 500
$ws.Cells.Item(581, 13).Value = 500
$ws.Cells.Item(581, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(581, 15).Value = "Región del Maule"
$ws.Cells.Item(581, 16).Value = 100
$ws.Cells.Item(581, 17).Value = 5
$ws.Cells.Item(581, 18).Value = "Hortaliza"
